$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 12, shifting existing rows 12..98 down to 13..99.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new weekly record.
$ws.Range("A12").Value = 4
$ws.Range("B12").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C12").Value = "Los Lagos"
$ws.Range("D12").Value = 45069
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = 100112043
$ws.Range("G12").Value = "Pepino dulce"
$ws.Range("H12").Value = "Cultivar IV Región"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 21000
$ws.Range("M12").Value = 20500
$ws.Range("N12").Value = "$/bandeja 18 kilos"
$ws.Range("O12").Value = "Provincia de Limarí"
$ws.Range("P12").Value = 1139
$ws.Range("Q12").Value = 18
$ws.Range("R12").Value = "Hortaliza"
